# Update Devguide's images (#134)
# - Refresh the cached "date last modified" figure-out field shown on
#   every slide layout + the slide master (6/7/2018 -> 24/10/18).
# - Rename "address book" -> "loan book" in the Undo/Redo activity
#   diagram on slide 1 (two shapes).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Date placeholder ("datetimeFigureOut" field) on every layout
#    and on the slide master: 6/7/2018 -> 24/10/18
# ---------------------------------------------------------------
$oldDate = "6/7/2018"
$newDate = "24/10/18"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if (-not $sh.HasTextFrame) { continue }
        $isDatePlaceholder = $false
        if ($sh.Type -eq 14) {
            # msoPlaceholder
            try {
                if ($sh.PlaceholderFormat.Type -eq 16) {
                    # ppPlaceholderDate
                    $isDatePlaceholder = $true
                }
            } catch {
                $isDatePlaceholder = $false
            }
        }
        if (-not $isDatePlaceholder -and $sh.Name -like "Date Placeholder*") {
            $isDatePlaceholder = $true
        }
        if ($isDatePlaceholder) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# ---------------------------------------------------------------
# 2) Slide 1 text: "address book" -> "loan book"
# ---------------------------------------------------------------
$slide = $p.Slides.Item(1)

# "[command commits address book]" -> "[command commits loan book]"
$shCommand = $slide.Shapes.Item(7)
$trCommand = $shCommand.TextFrame.TextRange
$fullCommand = $trCommand.Text
$oldCommandRun = "command commits address book]"
$newCommandRun = "command commits loan book]"
$startCommand = $fullCommand.IndexOf($oldCommandRun) + 1
if ($startCommand -gt 0) {
    $runRange = $trCommand.Characters($startCommand, $oldCommandRun.Length)
    $runRange.Text = $newCommandRun
}

# "Purge redundant states and then save address book to addressBookStateList "
# -> "Purge redundant states and then save loan book to loanBookStateList "
$shPurge = $slide.Shapes.Item(8)
$trPurge = $shPurge.TextFrame.TextRange
$fullPurge = $trPurge.Text

$oldLead = "Purge redundant states and then save address book to "
$newLead = "Purge redundant states and then save loan book to "
$startLead = $fullPurge.IndexOf($oldLead) + 1
if ($startLead -gt 0) {
    $leadRange = $trPurge.Characters($startLead, $oldLead.Length)
    $leadRange.Text = $newLead
}

# Re-read text/offsets after the first edit shifted character positions.
$fullPurge2 = $trPurge.Text
$oldVar = "addressBookStateList"
$newVar = "loanBookStateList"
$startVar = $fullPurge2.IndexOf($oldVar) + 1
if ($startVar -gt 0) {
    $varRange = $trPurge.Characters($startVar, $oldVar.Length)
    # Replacing straight across introduces a spurious leftover run because
    # the text-diff shares a common "...BookStateList" suffix with the old
    # value; clearing it with a disjoint placeholder first avoids that.
    $varRange.Text = "##################"
    $varRange2 = $trPurge.Characters($startVar, $oldVar.Length - 2)
    $varRange2.Text = $newVar
}
